$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'ECs'
$ws.Range("B2").Value = 'Hp'
$ws.Range("C2").Value = 'Itgb2'
$ws.Range("D2").Value = 'ECs'
$ws.Range("E2").Value = [double]"1"
$ws.Range("F2").Value = [double]"0.3333333333333333"
$ws.Range("G2").Value = [double]"0.1238986666666667"
$ws.Range("H2").Value = [double]"0.371696"
$ws.Range("I2").Value = [double]"0.01923905185495286"
$ws.Range("J2").Value = [double]"0.01923905185495286"
$ws.Range("K2").Value = [double]"1"
$ws.Range("L2").Value = [double]"0.3333333333333333"
$ws.Range("M2").Value = [double]"0.1145113333333333"
$ws.Range("N2").Value = [double]"0.343534"
$ws.Range("O2").Value = [double]"0.001785365609625045"
$ws.Range("P2").Value = [double]"0.001785365609625044"
$ws.Range("Q2").Value = [double]"0.01418780151822222"
$ws.Range("R2").Value = [double]"0.127690213664"
$ws.Range("S2").Value = [double]"3.434874154362576E-05"
$ws.Range("T2").Value = [double]"3.434874154362576E-05"

$ws.Range("A3").Value = 'ECs'
$ws.Range("B3").Value = 'Hp'
$ws.Range("C3").Value = 'Itgb2'
$ws.Range("D3").Value = 'FAPs'
$ws.Range("E3").Value = [double]"1"
$ws.Range("F3").Value = [double]"0.3333333333333333"
$ws.Range("G3").Value = [double]"0.1238986666666667"
$ws.Range("H3").Value = [double]"0.371696"
$ws.Range("I3").Value = [double]"0.01923905185495286"
$ws.Range("J3").Value = [double]"0.01923905185495286"
$ws.Range("K3").Value = [double]"3"
$ws.Range("L3").Value = [double]"1"
$ws.Range("M3").Value = [double]"0.467525"
$ws.Range("N3").Value = [double]"1.402575"
$ws.Range("O3").Value = [double]"0.007289261528465441"
$ws.Range("P3").Value = [double]"0.007289261528465441"
$ws.Range("Q3").Value = [double]"0.05792572413333334"
$ws.Range("R3").Value = [double]"0.5213315172"
$ws.Range("S3").Value = [double]"0.0001402384805304596"
$ws.Range("T3").Value = [double]"0.0001402384805304596"

$ws.Range("A4").Value = 'ECs'
$ws.Range("B4").Value = 'Hp'
$ws.Range("C4").Value = 'Itgb2'
$ws.Range("D4").Value = 'Resolving-Mac'
$ws.Range("E4").Value = [double]"1"
$ws.Range("F4").Value = [double]"0.3333333333333333"
$ws.Range("G4").Value = [double]"0.1238986666666667"
$ws.Range("H4").Value = [double]"0.371696"
$ws.Range("I4").Value = [double]"0.01923905185495286"
$ws.Range("J4").Value = [double]"0.01923905185495286"
$ws.Range("K4").Value = [double]"3"
$ws.Range("L4").Value = [double]"1"
$ws.Range("M4").Value = [double]"63.556834"
$ws.Range("N4").Value = [double]"190.670502"
$ws.Range("O4").Value = [double]"0.9909253728619096"
$ws.Range("P4").Value = [double]"0.9909253728619095"
$ws.Range("Q4").Value = [double]"7.874606990154668"
$ws.Range("R4").Value = [double]"70.871462911392"
$ws.Range("S4").Value = [double]"0.01906446463287878"
$ws.Range("T4").Value = [double]"0.01906446463287878"

$ws.Range("A5").Value = 'FAPs'
$ws.Range("B5").Value = 'Hp'
$ws.Range("C5").Value = 'Itgb2'
$ws.Range("D5").Value = 'ECs'
$ws.Range("E5").Value = [double]"3"
$ws.Range("F5").Value = [double]"1"
$ws.Range("G5").Value = [double]"5.810518333333333"
$ws.Range("H5").Value = [double]"17.431555"
$ws.Range("I5").Value = [double]"0.9022604239955847"
$ws.Range("J5").Value = [double]"0.9022604239955845"
$ws.Range("K5").Value = [double]"1"
$ws.Range("L5").Value = [double]"0.3333333333333333"
$ws.Range("M5").Value = [double]"0.1145113333333333"
$ws.Range("N5").Value = [double]"0.343534"
$ws.Range("O5").Value = [double]"0.001785365609625045"
$ws.Range("P5").Value = [double]"0.001785365609625044"
$ws.Range("Q5").Value = [double]"0.6653702017077778"
$ws.Range("R5").Value = [double]"5.98833181537"
$ws.Range("S5").Value = [double]"0.001610864731927428"
$ws.Range("T5").Value = [double]"0.001610864731927428"

$ws.Range("A6").Value = 'FAPs'
$ws.Range("B6").Value = 'Hp'
$ws.Range("C6").Value = 'Itgb2'
$ws.Range("D6").Value = 'FAPs'
$ws.Range("E6").Value = [double]"3"
$ws.Range("F6").Value = [double]"1"
$ws.Range("G6").Value = [double]"5.810518333333333"
$ws.Range("H6").Value = [double]"17.431555"
$ws.Range("I6").Value = [double]"0.9022604239955847"
$ws.Range("J6").Value = [double]"0.9022604239955845"
$ws.Range("K6").Value = [double]"3"
$ws.Range("L6").Value = [double]"1"
$ws.Range("M6").Value = [double]"0.467525"
$ws.Range("N6").Value = [double]"1.402575"
$ws.Range("O6").Value = [double]"0.007289261528465441"
$ws.Range("P6").Value = [double]"0.007289261528465441"
$ws.Range("Q6").Value = [double]"2.716562583791667"
$ws.Range("R6").Value = [double]"24.449063254125"
$ws.Range("S6").Value = [double]"0.006576812197287932"
$ws.Range("T6").Value = [double]"0.00657681219728793"

$ws.Range("A7").Value = 'FAPs'
$ws.Range("B7").Value = 'Hp'
$ws.Range("C7").Value = 'Itgb2'
$ws.Range("D7").Value = 'Resolving-Mac'
$ws.Range("E7").Value = [double]"3"
$ws.Range("F7").Value = [double]"1"
$ws.Range("G7").Value = [double]"5.810518333333333"
$ws.Range("H7").Value = [double]"17.431555"
$ws.Range("I7").Value = [double]"0.9022604239955847"
$ws.Range("J7").Value = [double]"0.9022604239955845"
$ws.Range("K7").Value = [double]"3"
$ws.Range("L7").Value = [double]"1"
$ws.Range("M7").Value = [double]"63.556834"
$ws.Range("N7").Value = [double]"190.670502"
$ws.Range("O7").Value = [double]"0.9909253728619096"
$ws.Range("P7").Value = [double]"0.9909253728619095"
$ws.Range("Q7").Value = [double]"369.2981491656233"
$ws.Range("R7").Value = [double]"3323.68334249061"
$ws.Range("S7").Value = [double]"0.8940727470663694"
$ws.Range("T7").Value = [double]"0.8940727470663691"

$ws.Range("A8").Value = 'MuSCs'
$ws.Range("B8").Value = 'Hp'
$ws.Range("C8").Value = 'Itgb2'
$ws.Range("D8").Value = 'ECs'
$ws.Range("E8").Value = [double]"1"
$ws.Range("F8").Value = [double]"0.3333333333333333"
$ws.Range("G8").Value = [double]"0.1062546666666667"
$ws.Range("H8").Value = [double]"0.318764"
$ws.Range("I8").Value = [double]"0.01649928200866351"
$ws.Range("J8").Value = [double]"0.01649928200866351"
$ws.Range("K8").Value = [double]"1"
$ws.Range("L8").Value = [double]"0.3333333333333333"
$ws.Range("M8").Value = [double]"0.1145113333333333"
$ws.Range("N8").Value = [double]"0.343534"
$ws.Range("O8").Value = [double]"0.001785365609625045"
$ws.Range("P8").Value = [double]"0.001785365609625044"
$ws.Range("Q8").Value = [double]"0.01216736355288889"
$ws.Range("R8").Value = [double]"0.109506271976"
$ws.Range("S8").Value = [double]"2.945725068177307E-05"
$ws.Range("T8").Value = [double]"2.945725068177305E-05"

$ws.Range("A9").Value = 'MuSCs'
$ws.Range("B9").Value = 'Hp'
$ws.Range("C9").Value = 'Itgb2'
$ws.Range("D9").Value = 'FAPs'
$ws.Range("E9").Value = [double]"1"
$ws.Range("F9").Value = [double]"0.3333333333333333"
$ws.Range("G9").Value = [double]"0.1062546666666667"
$ws.Range("H9").Value = [double]"0.318764"
$ws.Range("I9").Value = [double]"0.01649928200866351"
$ws.Range("J9").Value = [double]"0.01649928200866351"
$ws.Range("K9").Value = [double]"3"
$ws.Range("L9").Value = [double]"1"
$ws.Range("M9").Value = [double]"0.467525"
$ws.Range("N9").Value = [double]"1.402575"
$ws.Range("O9").Value = [double]"0.007289261528465441"
$ws.Range("P9").Value = [double]"0.007289261528465441"
$ws.Range("Q9").Value = [double]"0.04967671303333333"
$ws.Range("R9").Value = [double]"0.4470904173"
$ws.Range("S9").Value = [double]"0.000120267581593053"
$ws.Range("T9").Value = [double]"0.0001202675815930529"

$ws.Range("A10").Value = 'MuSCs'
$ws.Range("B10").Value = 'Hp'
$ws.Range("C10").Value = 'Itgb2'
$ws.Range("D10").Value = 'Resolving-Mac'
$ws.Range("E10").Value = [double]"1"
$ws.Range("F10").Value = [double]"0.3333333333333333"
$ws.Range("G10").Value = [double]"0.1062546666666667"
$ws.Range("H10").Value = [double]"0.318764"
$ws.Range("I10").Value = [double]"0.01649928200866351"
$ws.Range("J10").Value = [double]"0.01649928200866351"
$ws.Range("K10").Value = [double]"3"
$ws.Range("L10").Value = [double]"1"
$ws.Range("M10").Value = [double]"63.556834"
$ws.Range("N10").Value = [double]"190.670502"
$ws.Range("O10").Value = [double]"0.9909253728619096"
$ws.Range("P10").Value = [double]"0.9909253728619095"
$ws.Range("Q10").Value = [double]"6.753210211058667"
$ws.Range("R10").Value = [double]"60.778891899528"
$ws.Range("S10").Value = [double]"0.01634955717638869"
$ws.Range("T10").Value = [double]"0.01634955717638869"

$ws.Range("A11").Value = 'Resolving-Mac'
$ws.Range("B11").Value = 'Hp'
$ws.Range("C11").Value = 'Itgb2'
$ws.Range("D11").Value = 'ECs'
$ws.Range("E11").Value = [double]"2"
$ws.Range("F11").Value = [double]"0.6666666666666666"
$ws.Range("G11").Value = [double]"0.3992853333333333"
$ws.Range("H11").Value = [double]"1.197856"
$ws.Range("I11").Value = [double]"0.06200124214079897"
$ws.Range("J11").Value = [double]"0.06200124214079897"
$ws.Range("K11").Value = [double]"1"
$ws.Range("L11").Value = [double]"0.3333333333333333"
$ws.Range("M11").Value = [double]"0.1145113333333333"
$ws.Range("N11").Value = [double]"0.343534"
$ws.Range("O11").Value = [double]"0.001785365609625045"
$ws.Range("P11").Value = [double]"0.001785365609625044"
$ws.Range("Q11").Value = [double]"0.04572269590044445"
$ws.Range("R11").Value = [double]"0.411504263104"
$ws.Range("S11").Value = [double]"0.0001106948854722176"
$ws.Range("T11").Value = [double]"0.0001106948854722175"

$ws.Range("A12").Value = 'Resolving-Mac'
$ws.Range("B12").Value = 'Hp'
$ws.Range("C12").Value = 'Itgb2'
$ws.Range("D12").Value = 'FAPs'
$ws.Range("E12").Value = [double]"2"
$ws.Range("F12").Value = [double]"0.6666666666666666"
$ws.Range("G12").Value = [double]"0.3992853333333333"
$ws.Range("H12").Value = [double]"1.197856"
$ws.Range("I12").Value = [double]"0.06200124214079897"
$ws.Range("J12").Value = [double]"0.06200124214079897"
$ws.Range("K12").Value = [double]"3"
$ws.Range("L12").Value = [double]"1"
$ws.Range("M12").Value = [double]"0.467525"
$ws.Range("N12").Value = [double]"1.402575"
$ws.Range("O12").Value = [double]"0.007289261528465441"
$ws.Range("P12").Value = [double]"0.007289261528465441"
$ws.Range("Q12").Value = [double]"0.1866758754666667"
$ws.Range("R12").Value = [double]"1.6800828792"
$ws.Range("S12").Value = [double]"0.0004519432690539962"
$ws.Range("T12").Value = [double]"0.0004519432690539962"

$ws.Range("A13").Value = 'Resolving-Mac'
$ws.Range("B13").Value = 'Hp'
$ws.Range("C13").Value = 'Itgb2'
$ws.Range("D13").Value = 'Resolving-Mac'
$ws.Range("E13").Value = [double]"2"
$ws.Range("F13").Value = [double]"0.6666666666666666"
$ws.Range("G13").Value = [double]"0.3992853333333333"
$ws.Range("H13").Value = [double]"1.197856"
$ws.Range("I13").Value = [double]"0.06200124214079897"
$ws.Range("J13").Value = [double]"0.06200124214079897"
$ws.Range("K13").Value = [double]"3"
$ws.Range("L13").Value = [double]"1"
$ws.Range("M13").Value = [double]"63.556834"
$ws.Range("N13").Value = [double]"190.670502"
$ws.Range("O13").Value = [double]"0.9909253728619096"
$ws.Range("P13").Value = [double]"0.9909253728619095"
$ws.Range("Q13").Value = [double]"25.37731164930133"
$ws.Range("R13").Value = [double]"228.395804843712"
$ws.Range("S13").Value = [double]"0.06143860398627277"
$ws.Range("T13").Value = [double]"0.06143860398627275"
